$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores text values (e.g. thousand-separated
# "59.959.20" or plain "554.74"). Excel auto-detects numeric-looking
# strings and would silently convert them to real numbers, which would
# not match the original inline-string text cells, so force those
# specific cells to Text format before writing the new values.
$textRows = @(5,6,8,10,11,13,18,19,20,21,22,23,24,25,26,30,31,32,33,35,36,38,39,40,41,42,43,44,45,46,47,48,49,51)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = '@'
}

$ws.Range('D2').Value = '60.161.11'
$ws.Range('E2').Value = '  +3.71%  '
$ws.Range('D3').Value = '2.434.66'
$ws.Range('E3').Value = '  +3.86%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '556.38'
$ws.Range('E5').Value = '  +3.12%  '
$ws.Range('D6').Value = '139.01'
$ws.Range('E6').Value = '  +2.67%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').Value = '  +1.23%  '
$ws.Range('E9').Value = '  +5.54%  '
$ws.Range('D10').Value = '5.82'
$ws.Range('E10').Value = '  +4.85%  '
$ws.Range('D11').Value = '0.361'
$ws.Range('E11').Value = '  +2.38%  '
$ws.Range('E12').Value = '  -2.06%  '
$ws.Range('D13').Value = '24.91'
$ws.Range('E13').Value = '  +4.84%  '
$ws.Range('D14').Value = '2.867.83'
$ws.Range('E14').Value = '  +3.74%  '
$ws.Range('D15').Value = '60.021.50'
$ws.Range('E15').Value = '  +3.49%  '
$ws.Range('E16').Value = '  +5.12%  '
$ws.Range('D17').Value = '2.430.90'
$ws.Range('E17').Value = '  +3.69%  '
$ws.Range('D18').Value = '11.48'
$ws.Range('E18').Value = '  +7.58%  '
$ws.Range('D19').Value = '4.44'
$ws.Range('E19').Value = '  +4.18%  '
$ws.Range('D20').Value = '335.46'
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('D21').Value = '6.92'
$ws.Range('E21').Value = '  +2.02%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = '64.77'
$ws.Range('E23').Value = '  +3.33%  '
$ws.Range('D24').Value = '0.170'
$ws.Range('E24').Value = '  +2.03%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = '8.56'
$ws.Range('E25').Value = '  +0.99%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').Value = '0.0₃0793'
$ws.Range('E28').Value = '  +8.03%  '
$ws.Range('E29').Value = '  +3.79%  '
$ws.Range('D30').Value = '171.30'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').Value = '6.32'
$ws.Range('E31').Value = '  +3.08%  '
$ws.Range('D32').Value = '18.79'
$ws.Range('E32').Value = '  +1.85%  '
$ws.Range('D33').Value = '1.03'
$ws.Range('E33').Value = '  -0.89%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = '1.32'
$ws.Range('E35').Value = '  +5.50%  '
$ws.Range('D36').Value = '4.27'
$ws.Range('E36').Value = '  +1.27%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '1.65'
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('D39').Value = '40.15'
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('D40').Value = '0.418'
$ws.Range('E40').Value = '  +10.81%  '
$ws.Range('D41').Value = '318.02'
$ws.Range('E41').Value = '  +8.70%  '
$ws.Range('D42').Value = '3.74'
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').Value = '143.37'
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range('D44').Value = '0.0964'
$ws.Range('E44').Value = '  +1.67%  '
$ws.Range('D45').Value = '0.0525'
$ws.Range('E45').Value = '  +4.47%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '19.44'
$ws.Range('E46').Value = '  +1.50%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.574'
$ws.Range('E47').Value = '  +2.24%  '
$ws.Range('B48').Value = 'Polygon'
$ws.Range('C48').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D48').Value = '0.406'
$ws.Range('E48').Value = '  +5.39%  '
$ws.Range('D49').Value = '0.0227'
$ws.Range('E49').Value = '  +3.25%  '
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').Value = '1.63'
$ws.Range('E51').Value = '  +5.39%  '
